# Populate the "Plan de Actividades" header row and apply formatting.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header values - entered in this specific order so the shared-string
# table is built in the same sequence as the authored workbook.
$ws.Range("A1").Value = "Tarea"
$ws.Range("D1").Value = "Responsable a cargo"
$ws.Range("C1").Value = "Tiempo estimado de completitud"
$ws.Range("B1").Value = "Descripción de la tarea"
$ws.Range("E1").Value = "Fecha de entrega"

# Column widths to fit the header text
$ws.Columns.Item(2).ColumnWidth = 35.21875
$ws.Columns.Item(3).ColumnWidth = 45.21875
$ws.Columns.Item(4).ColumnWidth = 32.44140625
$ws.Columns.Item(5).ColumnWidth = 25.6640625

# Row height for the header row
$ws.Rows.Item(1).RowHeight = 15.6

# Formatting: bold Arial 12, filled background, thin border all around
$headerRange = $ws.Range("A1:E1")
$headerRange.Font.Bold = $true
$headerRange.Font.Size = 12
$headerRange.Font.Name = "Arial"
$headerRange.Interior.Pattern = 1
$headerRange.Interior.ThemeColor = 8
$headerRange.Interior.TintAndShade = 0.79998168889431442

$headerRange.Borders.Item(7).LineStyle = 1
$headerRange.Borders.Item(7).Weight = 2
$headerRange.Borders.Item(8).LineStyle = 1
$headerRange.Borders.Item(8).Weight = 2
$headerRange.Borders.Item(9).LineStyle = 1
$headerRange.Borders.Item(9).Weight = 2
$headerRange.Borders.Item(10).LineStyle = 1
$headerRange.Borders.Item(10).Weight = 2

$ws.Range("B2").Select()
